$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 should now hold the "LightGBM" data (previously on row 7).
$ws.Range("A3").Value = "LightGBM"
$ws.Range("B3").Value = 0.8881542699724518
$ws.Range("C3").Value = 0.8888128394179855
$ws.Range("D3").Value = 0.8881542699724518
$ws.Range("E3").Value = 0.8870059506730712

# Remove the now-obsolete rows (old Support Vector Classifier, CART,
# Random Forest, LightGBM, XGBoost rows 4-8).
$ws.Range("A4:E8").EntireRow.Delete()
